# week 7 Meeting.pptx - apply the authored edits:
#  1. Slide 1 subtitle: "Week 6 Meeting, 12/11/2021" -> "Week 7 Meeting, 12/11/2021"
#  2. Cached date-placeholder text on the slide master and every slide layout:
#     "11/15/2021" -> "12/11/2021" (PowerPoint re-stamps these cached field
#     values to the current date whenever the deck is re-saved).

$p = $ppt.ActivePresentation

# --- 1. Fix the subtitle text on slide 1 --------------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Runs().Count -ge 1 -and $tr.Runs(1).Text.TrimEnd() -eq "Week 6 Meeting, 12/11/2021") {
            # Update only the run text so the rest of the paragraph
            # structure (trailing empty paragraphs, run properties) is
            # left untouched.
            $tr.Runs(1).Text = "Week 7 Meeting, 12/11/2021"
        }
    }
}

# --- 2. Refresh the cached "11/15/2021" date field text -----------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text.TrimEnd() -eq "11/15/2021") {
                $tr.Text = "12/11/2021"
            }
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DatePlaceholders $p.SlideMaster.CustomLayouts.Item($li).Shapes
}
